# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (Murcott - Especial/Primera/Segunda, Provincia de
# Limarí, date 44491) at the top of the Mandarina data block, pushing the
# existing rows 189-210 down to 192-213.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 189 (old first data row of
# this block). Inserting three times at the same row index pushes everything
# down by three rows total.
$ws.Rows.Item(189).Insert()
$ws.Rows.Item(189).Insert()
$ws.Rows.Item(189).Insert()

# New row 189: Murcott / Especial
$ws.Cells.Item(189, 1).Value = 5
$ws.Cells.Item(189, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(189, 3).Value = "Maule"
$ws.Cells.Item(189, 4).Value = 44491
$ws.Cells.Item(189, 5).Value = 7
$ws.Cells.Item(189, 6).Value = "Fruta"
$ws.Cells.Item(189, 7).Value = 100102
$ws.Cells.Item(189, 8).Value = "Cítricos"
$ws.Cells.Item(189, 9).Value = 100102004
$ws.Cells.Item(189, 10).Value = "Mandarina"
$ws.Cells.Item(189, 11).Value = "Murcott"
$ws.Cells.Item(189, 12).Value = "Especial"
$ws.Cells.Item(189, 13).Value = 140
$ws.Cells.Item(189, 14).Value = 7000
$ws.Cells.Item(189, 15).Value = 7000
$ws.Cells.Item(189, 16).Value = 7000
$ws.Cells.Item(189, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(189, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(189, 19).Value = 389
$ws.Cells.Item(189, 20).Value = 18

# New row 190: Murcott / Primera
$ws.Cells.Item(190, 1).Value = 5
$ws.Cells.Item(190, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(190, 3).Value = "Maule"
$ws.Cells.Item(190, 4).Value = 44491
$ws.Cells.Item(190, 5).Value = 7
$ws.Cells.Item(190, 6).Value = "Fruta"
$ws.Cells.Item(190, 7).Value = 100102
$ws.Cells.Item(190, 8).Value = "Cítricos"
$ws.Cells.Item(190, 9).Value = 100102004
$ws.Cells.Item(190, 10).Value = "Mandarina"
$ws.Cells.Item(190, 11).Value = "Murcott"
$ws.Cells.Item(190, 12).Value = "Primera"
$ws.Cells.Item(190, 13).Value = 180
$ws.Cells.Item(190, 14).Value = 6000
$ws.Cells.Item(190, 15).Value = 6000
$ws.Cells.Item(190, 16).Value = 6000
$ws.Cells.Item(190, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(190, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(190, 19).Value = 333
$ws.Cells.Item(190, 20).Value = 18

# New row 191: Murcott / Segunda
$ws.Cells.Item(191, 1).Value = 5
$ws.Cells.Item(191, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(191, 3).Value = "Maule"
$ws.Cells.Item(191, 4).Value = 44491
$ws.Cells.Item(191, 5).Value = 7
$ws.Cells.Item(191, 6).Value = "Fruta"
$ws.Cells.Item(191, 7).Value = 100102
$ws.Cells.Item(191, 8).Value = "Cítricos"
$ws.Cells.Item(191, 9).Value = 100102004
$ws.Cells.Item(191, 10).Value = "Mandarina"
$ws.Cells.Item(191, 11).Value = "Murcott"
$ws.Cells.Item(191, 12).Value = "Segunda"
$ws.Cells.Item(191, 13).Value = 100
$ws.Cells.Item(191, 14).Value = 4000
$ws.Cells.Item(191, 15).Value = 4000
$ws.Cells.Item(191, 16).Value = 4000
$ws.Cells.Item(191, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(191, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(191, 19).Value = 222
$ws.Cells.Item(191, 20).Value = 18

Write-Host "Inserted 3 new Mandarina rows (Murcott Especial/Primera/Segunda, 44491)"
